$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (order matters for shared-string table layout) ---
$ws.Range("A1").Value2 = "Tree"
$ws.Range("C1").Value2 = "Dendrometer"
$ws.Range("B1").Value2 = "Letter"

# --- Row 9: drop the stray "applied fill" flag left over from earlier editing ---
# (the fill itself was never actually colored - this just removes the dead
#  cellXfs entries / customFormat flag so the row matches the normal rows)
$ws.Rows.Item(9).ClearFormats() | Out-Null
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("F9").ClearContents() | Out-Null

# --- Column widths: re-fit after the header text changed width ---
$ws.Columns.Item(1).ColumnWidth = 3.830729167
$ws.Columns.Item(2).ColumnWidth = 5.053385417
$ws.Columns.Item(3).ColumnWidth = 11.05338542
$ws.Columns.Item(5).ColumnWidth = 7.721354167
$ws.Columns.Item(6).ColumnWidth = 9.944010417
$ws.Columns.Item(8).ColumnWidth = 34.83072917

# --- Selection moved by the author before saving ---
$ws.Range("D4").Select() | Out-Null
